$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Text content change: the "lambda in sequence provided by glmnet::glmnet
#    corresponding to maximum F1 measure" label is replaced with a new,
#    shorter label "lambda in sequence provided by glmnet::cv.glmnet".
# ---------------------------------------------------------------------------
$ws.Range("C5").Value = "lambda in sequence provided by glmnet::cv.glmnet"

# ---------------------------------------------------------------------------
# 2. Updated data values (re-run of the underlying analysis produced new
#    numbers across most of the data rows).
# ---------------------------------------------------------------------------
$ws.Range("E3").Value = 0.18
$ws.Range("F3").Value = 0.69

$ws.Range("E5").Value = 0.0001

$ws.Range("H6").Value = 1

$ws.Range("D7").Value = 0.9844
$ws.Range("E7").Value = 0.9811
$ws.Range("F7").Value = 0.903
$ws.Range("G7").Value = 0.9656
$ws.Range("H7").Value = 0.9479

$ws.Range("D8").Value = 0.9993
$ws.Range("E8").Value = 0.9992
$ws.Range("F8").Value = 0.9494
$ws.Range("G8").Value = 0.9983
$ws.Range("H8").Value = 0.9772

$ws.Range("D9").Value = 0.9962
$ws.Range("E9").Value = 0.9957
$ws.Range("F9").Value = 0.9937
$ws.Range("G9").Value = 0.9961
$ws.Range("H9").Value = 0.9965

$ws.Range("D10").Value = 0.9667
$ws.Range("E10").Value = 0.9444
$ws.Range("F10").Value = 0.8097
$ws.Range("G10").Value = 0.9369
$ws.Range("H10").Value = 0.9457

$ws.Range("D11").Value = 0.0029
$ws.Range("E11").Value = 0.0026
$ws.Range("F11").Value = 0.0002
$ws.Range("G11").Value = 0.0019
$ws.Range("H11").Value = 0.0018

$ws.Range("D12").Value = 0.9179
$ws.Range("E12").Value = 0.9232
$ws.Range("F12").Value = 0.9916
$ws.Range("G12").Value = 0.9412
$ws.Range("H12").Value = 0.9461

$ws.Range("D13").Value = 0.9412
$ws.Range("E13").Value = 0.9335
$ws.Range("F13").Value = 0.8912
$ws.Range("G13").Value = 0.9389
$ws.Range("H13").Value = 0.9458

# ---------------------------------------------------------------------------
# 3. Number formats: the percentage-ish "accounting" style format now shows
#    four decimal places instead of three, and the "optimal K" row (H6, an
#    integer count) switches to a plain thousands-style integer format.
# ---------------------------------------------------------------------------
$numFmt4Dec = "_(* #,##0.0000_);_(* \(#,##0.0000\);_(* ""-""????_);_(@_)"
$numFmtInt  = "_(* #,##0_);_(* \(#,##0\);_(* ""-""_);_(@_)"

$ws.Range("D2:H13").NumberFormat = $numFmt4Dec
$ws.Range("H6").NumberFormat = $numFmtInt

# ---------------------------------------------------------------------------
# 4. Header row (B2:H2) becomes centered.
# ---------------------------------------------------------------------------
$ws.Range("B2:H2").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# 5. Highlight-color reshuffle within the "average" rows (7-13): a handful of
#    cells swap their fill colors.
# ---------------------------------------------------------------------------
$red    = 255
$orange = 49407
$yellow = 65535
$green  = 5296274
$blue   = 15773696

$ws.Range("E8").Interior.Color = $orange
$ws.Range("G8").Interior.Color = $yellow

$ws.Range("D9").Interior.Color = $orange
$ws.Range("H9").Interior.Color = $red

$ws.Range("E11").Interior.Color = $green
$ws.Range("F11").Interior.Color = $red
$ws.Range("G11").Interior.Color = $yellow
$ws.Range("H11").Interior.Color = $orange

$ws.Range("D13").Interior.Color = $orange
$ws.Range("H13").Interior.Color = $red

# ---------------------------------------------------------------------------
# 6. Column width adjustments (cosmetic re-fit of columns B-H).
# ---------------------------------------------------------------------------
$ws.Range("B1").EntireColumn.ColumnWidth = 41.329
$ws.Range("C1").EntireColumn.ColumnWidth = 66.329
$ws.Range("D1").EntireColumn.ColumnWidth = 18.16
$ws.Range("E1").EntireColumn.ColumnWidth = 23.661
$ws.Range("F1:H1").EntireColumn.ColumnWidth = 6.997
